$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "29.903.73"
$ws.Range("D3").Value = "1.888.62"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'0.7655"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").Value = "'242.84"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.3125"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "'25.67"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'0.07181"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "'0.08514"
$ws.Range("E11").Value = "  +4.67%  "
$ws.Range("D12").Value = "'0.7639"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "'5.362"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").Value = "1.870.29"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'93.61"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "'6.142"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "29.777.74"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "'244.39"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'0.000007806"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9988"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.119.60"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'7.980"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'0.1621"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "'9.419"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "'162.13"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").Value = "'2.035"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'1.466"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").Value = "'1.531"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "'4.497"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "'4.097"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "'0.05456"
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").Value = "'1.245"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "'0.7414"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("D37").Value = "'0.9999"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'0.4466"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "1.103.42"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("D43").Value = "'73.03"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "'6.068"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").Value = "'0.8523"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'102.81"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("D49").Value = "'7.657"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").Value = "'2.988"
$ws.Range("E50").Value = "  -4.04%  "
$ws.Range("D51").Value = "2.016.72"
$ws.Range("E51").Value = "  -1.28%  "
